$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append after the existing data (row 269, date serial 44343 /
# 2021-05-27), extending coverage through 2021-06-28 (commit: "aggiornamento
# fino a 28/06 incluso"). Columns: A = date serial, B = nuovi pos.,
# C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti.
$newRows = @(
    @(44344, 0, 0, 0),
    @(44345, 1, 1, 16.63616702711695),
    @(44346, 0, 1, 16.63616702711695),
    @(44347, 2, 3, 49.90850108135086),
    @(44348, 0, 3, 49.90850108135086),
    @(44349, 0, 3, 49.90850108135086),
    @(44350, 0, 3, 49.90850108135086),
    @(44351, 0, 3, 49.90850108135086),
    @(44352, 1, 3, 49.90850108135086),
    @(44353, 0, 3, 49.90850108135086),
    @(44354, 0, 1, 16.63616702711695),
    @(44355, 0, 1, 16.63616702711695),
    @(44356, 0, 1, 16.63616702711695),
    @(44357, 0, 1, 16.63616702711695),
    @(44358, 0, 1, 16.63616702711695),
    @(44359, 0, 0, 0),
    @(44360, 1, 1, 16.63616702711695),
    @(44361, 0, 1, 16.63616702711695),
    @(44362, 0, 1, 16.63616702711695),
    @(44363, 0, 1, 16.63616702711695),
    @(44364, 0, 1, 16.63616702711695),
    @(44365, 0, 1, 16.63616702711695),
    @(44366, 0, 1, 16.63616702711695),
    @(44367, 0, 0, 0),
    @(44368, 0, 0, 0),
    @(44369, 0, 0, 0),
    @(44370, 0, 0, 0),
    @(44371, 0, 0, 0),
    @(44372, 0, 0, 0),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$lastExistingRow = 269
$startRow = $lastExistingRow + 1
$endRow = $startRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Match the formatting of the last existing data row (date style on column A,
# plain numbers on B:D) across the whole newly added block.
$ws.Range("A$lastExistingRow`:D$lastExistingRow").Copy()
$ws.Range("A$startRow`:D$endRow").PasteSpecial(-4122)  # xlPasteFormats
